$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting/style from H1 (bold, bordered, centered) onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the data for I2:J14
$data = @(
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(11, 11),
    @(6, 7),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(5, 6),
    @(6, 7),
    @(7, 7),
    @(7, 7),
    @(9, 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
